# Refresh the cryptos worksheet with the latest price (column D) and
# 1-hour volume-change percentage (column E) figures pulled by the
# scheduled GitHub Actions job.
#
# Price/percentage values are stored as plain text in this sheet, so a
# leading apostrophe is used for numeric-looking prices (e.g. '7.70)
# to stop Excel from auto-converting them to numbers and silently
# dropping meaningful trailing zeros.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.629.23"
$ws.Range("E2").Value = "  -2.09%  "

$ws.Range("D3").Value = "2.344.42"
$ws.Range("E3").Value = "  -3.00%  "

$ws.Range("E4").Value = "  -0.18%  "

$ws.Range("D5").Value = "'320.24"
$ws.Range("E5").Value = "  -1.46%  "

$ws.Range("D6").Value = "'105.78"
$ws.Range("E6").Value = "  +0.85%  "

$ws.Range("D7").Value = "'0.637"
$ws.Range("E7").Value = "  -1.42%  "

$ws.Range("E8").Value = "  +0.05%  "

$ws.Range("E9").Value = "  -7.20%  "

$ws.Range("D10").Value = "'40.85"
$ws.Range("E10").Value = "  -3.49%  "

$ws.Range("E11").Value = "  -3.12%  "

$ws.Range("E12").Value = "  -3.50%  "

$ws.Range("E13").Value = "  -5.35%  "

$ws.Range("E14").Value = "  -0.60%  "

$ws.Range("D15").Value = "'15.95"
$ws.Range("E15").Value = "  -7.31%  "

$ws.Range("D16").Value = "2.696.53"
$ws.Range("E16").Value = "  -3.14%  "

$ws.Range("D17").Value = "2.350.84"
$ws.Range("E17").Value = "  -2.87%  "

$ws.Range("D18").Value = "42.593.18"
$ws.Range("E18").Value = "  -2.25%  "

$ws.Range("D19").Value = "'7.70"
$ws.Range("E19").Value = "  +3.01%  "

$ws.Range("E20").Value = "  -4.32%  "

$ws.Range("D21").Value = "'77.20"
$ws.Range("E21").Value = "  +2.07%  "

$ws.Range("D22").Value = "'3.56"
$ws.Range("E22").Value = "  +1.73%  "

$ws.Range("D23").Value = "'258.96"
$ws.Range("E23").Value = "  -0.66%  "

$ws.Range("D24").Value = "'2.31"
$ws.Range("E24").Value = "  -5.44%  "

$ws.Range("D25").Value = "'9.54"
$ws.Range("E25").Value = "  -1.11%  "

$ws.Range("E26").Value = "  +0.08%  "

$ws.Range("D27").Value = "'11.40"
$ws.Range("E27").Value = "  -4.89%  "

$ws.Range("D28").Value = "'23.13"
$ws.Range("E28").Value = "  +0.88%  "

$ws.Range("E29").Value = "  -0.98%  "

$ws.Range("D30").Value = "'174.29"
$ws.Range("E30").Value = "  -3.21%  "

$ws.Range("D31").Value = "'36.19"
$ws.Range("E31").Value = "  -5.55%  "

$ws.Range("D32").Value = "'3.01"
$ws.Range("E32").Value = "  -6.84%  "

$ws.Range("E33").Value = "  -4.88%  "

$ws.Range("D34").Value = "'6.05"
$ws.Range("E34").Value = "  +1.90%  "

$ws.Range("E35").Value = "  -1.93%  "

$ws.Range("E36").Value = "  +4.83%  "

$ws.Range("E37").Value = "  -6.07%  "

$ws.Range("E38").Value = "  -4.61%  "

$ws.Range("D39").Value = "'3.79"
$ws.Range("E39").Value = "  -5.66%  "

$ws.Range("E40").Value = "  -8.23%  "

$ws.Range("E41").Value = "  -10.92%  "

$ws.Range("D42").Value = "'70.27"
$ws.Range("E42").Value = "  +0.67%  "

$ws.Range("E43").Value = "  -1.76%  "

$ws.Range("E44").Value = "  -0.12%  "

$ws.Range("D45").Value = "'114.06"
$ws.Range("E45").Value = "  -9.32%  "

$ws.Range("D46").Value = "'11.87"
$ws.Range("E46").Value = "  -5.97%  "

$ws.Range("D47").Value = "'5.49"
$ws.Range("E47").Value = "  -3.59%  "

$ws.Range("D48").Value = "'9.18"
$ws.Range("E48").Value = "  -3.85%  "

$ws.Range("D49").Value = "'82.99"
$ws.Range("E49").Value = "  +5.88%  "

$ws.Range("D50").Value = "'72.85"
$ws.Range("E50").Value = "  -2.95%  "
